$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: merge the split "Date: 12/" + "9" + "/21" runs (under the
# 12/9/21 day header) into a single run "Date: 12/9/21", keeping the
# bold/size-36 run formatting and the <w:lastRenderedPageBreak/> marker.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$dateParaIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Date: 12/9/21`r") {
        $dateParaIndex = $i
        break
    }
}

if ($dateParaIndex -ne -1) {
    $dateRng = $d.Paragraphs.Item($dateParaIndex).Range
    # Range covering the paragraph's runs only (exclude the paragraph mark)
    $dateTarget = $d.Range($dateRng.Start, $dateRng.End - 1)
    $dateXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:lastRenderedPageBreak/><w:t>Date: 12/9/21</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $dateTarget.InsertXML($dateXml)
}

# ---------------------------------------------------------------------------
# Edit 2: add "Barplot and pie chart " (bold, with spell-check proofErr
# markers around "Barplot") as the content of the empty "Done:" list item
# under "Member Name: Megan " in the 12/9/21 section.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$barplotParaIndex = -1
for ($i = $dateParaIndex; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r" -and $i -gt 2) {
        $prevText = $d.Paragraphs.Item($i - 1).Range.Text
        if ($prevText -eq "Done: `r") {
            $memberText = $d.Paragraphs.Item($i - 2).Range.Text
            if ($memberText -like "*Megan*") {
                $barplotParaIndex = $i
                break
            }
        }
    }
}

if ($barplotParaIndex -ne -1) {
    $barplotRng = $d.Paragraphs.Item($barplotParaIndex).Range
    # Collapsed point at the very start of the (empty) paragraph's content,
    # so InsertXML fills the paragraph body without touching its <w:pPr>.
    $barplotTarget = $d.Range($barplotRng.Start, $barplotRng.Start)
    $barplotXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Barplot</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> and pie chart </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $barplotTarget.InsertXML($barplotXml)
}
